# Corresponding edits to High GDP files
# On the "VTStFES" sheet, row 6 (B6:H6) should now be formulas that mirror
# row 2 (B2:H2) instead of hard-coded values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VTStFES")

$ws.Range("B6").Formula = "=B2"
$ws.Range("C6:H6").Formula = "=C2"

# Set the active selection to match the edited range (cosmetic, matches diff)
$ws.Range("B6:H6").Select()

$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("K6").Select()
